$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-08 Saturday" "2025-11-09 Sunday"

Replace-Text "53÷5=10, 3" "55÷3=18, 1"
Replace-Text "82÷5=16, 2" "66÷2=33, 0"
Replace-Text "80÷7=11, 3" "75÷8=9, 3"
Replace-Text "19÷5=3, 4" "14÷8=1, 6"
Replace-Text "50÷5=10, 0" "35÷6=5, 5"

Replace-Text "71÷7=10, 1" "55÷8=6, 7"
Replace-Text "80÷9=8, 8" "77÷9=8, 5"
Replace-Text "86÷9=9, 5" "61÷7=8, 5"
Replace-Text "79÷4=19, 3" "99÷8=12, 3"
Replace-Text "76÷7=10, 6" "32÷5=6, 2"

Replace-Text "29÷6=4, 5" "77÷5=15, 2"
Replace-Text "57÷3=19, 0" "97÷8=12, 1"
Replace-Text "61÷8=7, 5" "47÷7=6, 5"
Replace-Text "16÷6=2, 4" "93÷8=11, 5"
Replace-Text "17÷2=8, 1" "64÷8=8, 0"

Replace-Text "44÷9=4, 8" "85÷2=42, 1"
Replace-Text "11÷4=2, 3" "75÷2=37, 1"
Replace-Text "79÷6=13, 1" "77÷3=25, 2"
Replace-Text "41÷3=13, 2" "43÷8=5, 3"
Replace-Text "46÷9=5, 1" "18÷7=2, 4"

Replace-Text "18÷2=9, 0" "31÷9=3, 4"
Replace-Text "57÷7=8, 1" "11÷3=3, 2"
Replace-Text "23÷7=3, 2" "54÷9=6, 0"
Replace-Text "90÷2=45, 0" "71÷6=11, 5"
Replace-Text "86÷7=12, 2" "35÷7=5, 0"
